{"js": "// Replace the date line and the multiplication problems in the table\n// with their new values, per the commit's diff. Each old text is\n// unique in the document, so a simple search+replace per pair is safe.\nconst replacements = [\n  [\"2025-12-13 Saturday\", \"2025-12-14 Sunday\"],\n  [\"141\u00d76=\", \"444\u00d76=\"],\n  [\"578\u00d73=\", \"656\u00d78=\"],\n  [\"471\u00d78=\", \"116\u00d79=\"],\n  [\"419\u00d77=\", \"967\u00d79=\"],\n  [\"840\u00d79=\", \"678\u00d76=\"],\n  [\"603\u00d77=\", \"279\u00d76=\"],\n  [\"423\u00d75=\", \"535\u00d73=\"],\n  [\"476\u00d78=\", \"833\u00d77=\"],\n  [\"125\u00d72=\", \"542\u00d73=\"],\n  [\"533\u00d77=\", \"186\u00d73=\"],\n  [\"282\u00d76=\", \"134\u00d74=\"],\n  [\"764\u00d75=\", \"749\u00d76=\"],\n  [\"787\u00d77=\", \"895\u00d75=\"],\n  [\"552\u00d77=\", \"785\u00d76=\"],\n  [\"345\u00d78=\", \"696\u00d77=\"],\n  [\"454\u00d78=\", \"368\u00d74=\"],\n  [\"298\u00d79=\", \"942\u00d76=\"],\n  [\"382\u00d73=\", \"880\u00d78=\"],\n  [\"406\u00d79=\", \"251\u00d76=\"],\n  [\"304\u00d73=\", \"911\u00d72=\"],\n  [\"858\u00d72=\", \"793\u00d75=\"],\n  [\"588\u00d78=\", \"690\u00d74=\"],\n  [\"395\u00d75=\", \"114\u00d76=\"],\n  [\"770\u00d79=\", \"501\u00d74=\"],\n  [\"366\u00d79=\", \"782\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the multiplication problems in the table\n# with their new values, per the commit's diff. Each old text is\n# unique in the document, so Find/Replace (replace-all) per pair is safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-12-13 Saturday\", \"2025-12-14 Sunday\"),\n    @(\"141\u00d76=\", \"444\u00d76=\"),\n    @(\"578\u00d73=\", \"656\u00d78=\"),\n    @(\"471\u00d78=\", \"116\u00d79=\"),\n    @(\"419\u00d77=\", \"967\u00d79=\"),\n    @(\"840\u00d79=\", \"678\u00d76=\"),\n    @(\"603\u00d77=\", \"279\u00d76=\"),\n    @(\"423\u00d75=\", \"535\u00d73=\"),\n    @(\"476\u00d78=\", \"833\u00d77=\"),\n    @(\"125\u00d72=\", \"542\u00d73=\"),\n    @(\"533\u00d77=\", \"186\u00d73=\"),\n    @(\"282\u00d76=\", \"134\u00d74=\"),\n    @(\"764\u00d75=\", \"749\u00d76=\"),\n    @(\"787\u00d77=\", \"895\u00d75=\"),\n    @(\"552\u00d77=\", \"785\u00d76=\"),\n    @(\"345\u00d78=\", \"696\u00d77=\"),\n    @(\"454\u00d78=\", \"368\u00d74=\"),\n    @(\"298\u00d79=\", \"942\u00d76=\"),\n    @(\"382\u00d73=\", \"880\u00d78=\"),\n    @(\"406\u00d79=\", \"251\u00d76=\"),\n    @(\"304\u00d73=\", \"911\u00d72=\"),\n    @(\"858\u00d72=\", \"793\u00d75=\"),\n    @(\"588\u00d78=\", \"690\u00d74=\"),\n    @(\"395\u00d75=\", \"114\u00d76=\"),\n    @(\"770\u00d79=\", \"501\u00d74=\"),\n    @(\"366\u00d79=\", \"782\u00d78=\")\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll) | Out-Null\n}\n"}
